$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that contain the "sleep" entries (every 3rd row starting at row 2)
# whose N (JKL) and O (OS) columns flip from FALSE to TRUE.
$sleepRows = @(2,5,8,11,14,17,20,23,26,29,32,35,38,41,44,47,50,53,56,59,62,65,68,71,74,77)

foreach ($r in $sleepRows) {
    $ws.Range("N$r").Value = $true
    $ws.Range("O$r").Value = $true
}

# Row 78 ("activity" entry) - column I (StayStrong) flips from FALSE to TRUE.
$ws.Range("I78").Value = $true

# Append three new data rows for 2025-02-27 (sleep, activity, weekly_activity).
# Force the Date column to be stored as text (matching the existing date
# cells), then restore the default "Normal" style so no stray number
# format is left applied to the cells.
$dateCells = $ws.Range("A80:A82")
$dateCells.NumberFormat = "@"

$ws.Range("A80").Value = "2025-02-27"
$ws.Range("B80").Value = "sleep"
$ws.Range("C80").Value = $false
$ws.Range("D80").Value = $false
$ws.Range("E80").Value = $false
$ws.Range("F80").Value = $false
$ws.Range("G80").Value = $true
$ws.Range("H80").Value = $true
$ws.Range("I80").Value = $true
$ws.Range("J80").Value = $true
$ws.Range("K80").Value = $true
$ws.Range("L80").Value = $true
$ws.Range("M80").Value = $true
$ws.Range("N80").Value = $true
$ws.Range("O80").Value = $true

$ws.Range("A81").Value = "2025-02-27"
$ws.Range("B81").Value = "activity"
$ws.Range("C81").Value = $true
$ws.Range("D81").Value = $false
$ws.Range("E81").Value = $true
$ws.Range("F81").Value = $false
$ws.Range("G81").Value = $true
$ws.Range("H81").Value = $true
$ws.Range("I81").Value = $true
$ws.Range("J81").Value = $true
$ws.Range("K81").Value = $false
$ws.Range("L81").Value = $false
$ws.Range("M81").Value = $true
$ws.Range("N81").Value = $false
$ws.Range("O81").Value = $false

$ws.Range("A82").Value = "2025-02-27"
$ws.Range("B82").Value = "weekly_activity"
$ws.Range("C82").Value = $true
$ws.Range("D82").Value = $false
$ws.Range("E82").Value = $true
$ws.Range("F82").Value = $true
$ws.Range("G82").Value = $true
$ws.Range("H82").Value = $true
$ws.Range("I82").Value = $true
$ws.Range("J82").Value = $true
$ws.Range("K82").Value = $false
$ws.Range("L82").Value = $false
$ws.Range("M82").Value = $true
$ws.Range("N82").Value = $true
$ws.Range("O82").Value = $false

$dateCells.Style = "Normal"

Write-Host ("New dimension: " + $ws.UsedRange.Address())
